# Oscillazioni forzate smorzate: rename Sheet1, extend the data table with
# the forced-oscillation measurements (rows 2-17), and extend the trailing
# blank/styled rows through row 22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "oscillazioni forzate"

# --- Extend formatting into the new rows before writing values --------
# Row 11 (A:F) already carries the data-row style (s=6); clone it onto the
# rows that are brand-new or were previously blank / F-only (12-17).
$ws.Range("A11:F11").Copy()
$ws.Range("A12:F17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The old "omegaf" column, E11, used a one-off style (numFmtId 2, same look
# as s=6 but its own xf entry); collapse it onto the shared s=6 style used
# by the rest of column E.
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)

# Column F keeps its styled-but-empty trailer further down (rows 18-22).
$ws.Range("F18").Copy()
$ws.Range("F19:F22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data table (T, V, A, omega, omegaf) -------------------------------
$ws.Range("A2").Value = 3.91
$ws.Range("B2").Value = 1.9
$ws.Range("C2").Value = 2.95
$ws.Range("D2").Value = 1.61
$ws.Range("E2").Formula = "=6.28/A2"

$ws.Range("A3").Value = 3.55
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 3.43
$ws.Range("D3").Value = 1.78
$ws.Range("E3").Formula = "=2*PI()/A3"

$ws.Range("A4").Value = 3.2
$ws.Range("B4").Value = 2.1
$ws.Range("C4").Value = 4.35
$ws.Range("D4").Value = 2.03
$ws.Range("E4").Formula = "=2*PI()/A4"

$ws.Range("A5").Value = 3.04
$ws.Range("B5").Value = 2.2
$ws.Range("C5").Value = 5.5
$ws.Range("D5").Value = 2.15
$ws.Range("E5").Formula = "=2*PI()/A5"

$ws.Range("A6").Value = 2.79
$ws.Range("B6").Value = 2.3
$ws.Range("C6").Value = 6.8
$ws.Range("D6").Value = 2.23
$ws.Range("E6").Formula = "=2*PI()/A6"

$ws.Range("A7").Value = 2.66
$ws.Range("B7").Value = 2.35
$ws.Range("C7").Value = 9.39
$ws.Range("D7").Value = 2.37
$ws.Range("E7").Formula = "=2*PI()/A7"

$ws.Range("A8").Value = 2.54
$ws.Range("B8").Value = 2.4
$ws.Range("C8").Value = 11.9
$ws.Range("D8").Value = 2.49
$ws.Range("E8").Formula = "=6.28/A8"

$ws.Range("A9").Value = 2.49
$ws.Range("B9").Value = 2.45
$ws.Range("C9").Value = 11.8
$ws.Range("D9").Value = 2.5
$ws.Range("E9").Formula = "=2*PI()/A9"

$ws.Range("A10").Value = 2.4
$ws.Range("B10").Value = 2.5
$ws.Range("C10").Value = 11.8
$ws.Range("D10").Value = 2.66
$ws.Range("E10").Formula = "=2*PI()/A10"

$ws.Range("A11").Value = 2.35
$ws.Range("B11").Value = 2.6
$ws.Range("C11").Value = 10.5
$ws.Range("D11").Value = 2.71
$ws.Range("E11").Formula = "=2*PI()/A11"

$ws.Range("A12").Value = 2.31
$ws.Range("B12").Value = 2.65
$ws.Range("C12").Value = 9.98
$ws.Range("D12").Value = 2.73
$ws.Range("E12").Formula = "=2*PI()/A12"

$ws.Range("A13").Value = 2.18
$ws.Range("B13").Value = 2.7
$ws.Range("C13").Value = 6.28
$ws.Range("D13").Value = 2.86
$ws.Range("E13").Formula = "=2*PI()/A13"

$ws.Range("A14").Value = 2.07
$ws.Range("B14").Value = 2.8
$ws.Range("C14").Value = 4.85
$ws.Range("D14").Value = 3.06
$ws.Range("E14").Formula = "=2*PI()/A14"

$ws.Range("A15").Value = 2
$ws.Range("B15").Value = 2.9
$ws.Range("C15").Value = 3.75
$ws.Range("D15").Value = 3.14
$ws.Range("E15").Formula = "=2*PI()/A15"

$ws.Range("A16").Value = 1.92
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = 3.5
$ws.Range("D16").Value = 3.29
$ws.Range("E16").Formula = "=2*PI()/A16"

$ws.Range("A17").Value = 1.81
$ws.Range("B17").Value = 3.1
$ws.Range("C17").Value = 2.9
$ws.Range("D17").Value = 3.55
$ws.Range("E17").Formula = "=2*PI()/A17"
